$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.281.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "'3.131.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'602.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("D6").Value = "'142.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.37%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'3.127.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  +2.92%  "

$ws.Range("D12").Value = "'0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("E13").Value = "  +2.79%  "

$ws.Range("D14").Value = "'35.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'3.646.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "

$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").Value = "'64.193.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "'3.129.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "

$ws.Range("D19").Value = "'6.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "

$ws.Range("D20").Value = "'480.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").Value = "'14.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").Value = "'0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.68%  "

$ws.Range("D23").Value = "'7.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("D24").Value = "'85.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").Value = "'13.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").Value = "'8.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").Value = "'7.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.88%  "

$ws.Range("E30").Value = "  -4.15%  "

$ws.Range("D31").Value = "'0.113"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("D33").Value = "'26.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").Value = "'2.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.68%  "

$ws.Range("E35").Value = "  -0.78%  "

$ws.Range("D36").Value = "'0.0₃0770"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.54%  "

$ws.Range("D37").Value = "'5.96"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("D39").Value = "'3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.39%  "

$ws.Range("D40").Value = "'445.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.29%  "

$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("E43").Value = "  -1.50%  "

$ws.Range("D44").Value = "'2.855.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").Value = "'0.260"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.15%  "

$ws.Range("E46").Value = "  -1.64%  "

$ws.Range("D47").Value = "'2.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("D49").Value = "'26.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").Value = "'120.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.83%  "
